# Revise responsive design implementation
# Simplify code structure to improve readability.
#
# Appends a new row 68 (mirroring the structure of the existing rows) to
# each of the four worksheets, extending the used range from A1:I67 to
# A1:I68.

$wb = $excel.ActiveWorkbook

# time value (days since 1899-12-30) shared by every new row, formatted as
# the same custom date/time display the rest of column A already uses.
$newTime = 45854.46606481481
$dateFmt = "YYYY-MM-DD HH:MM:SS"

function Add-DataRow($Sheet, $Row, $BVal, $CVal, $DVal, $EVal, $FVal, $GVal, $HVal, $IVal) {
    $Sheet.Cells.Item($Row, 1).Value2 = $newTime
    $Sheet.Cells.Item($Row, 1).NumberFormat = $dateFmt

    $Sheet.Cells.Item($Row, 2).Value2 = $BVal
    $Sheet.Cells.Item($Row, 3).Value2 = $CVal
    $Sheet.Cells.Item($Row, 4).Value2 = $DVal
    $Sheet.Cells.Item($Row, 5).Value2 = $EVal

    $Sheet.Cells.Item($Row, 6).Value2 = $FVal
    $Sheet.Cells.Item($Row, 7).Value2 = $GVal
    $Sheet.Cells.Item($Row, 8).Value2 = $HVal
    $Sheet.Cells.Item($Row, 9).Value2 = $IVal
}

# --- MID_LFT_#1 ---
$ws1 = $wb.Worksheets.Item(1)
$g1 = [double]"5.68631262647113e+23"
Add-DataRow $ws1 68 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," `
    "0x01,0x50" `
    "0x07" `
    400 `
    $g1 `
    336 `
    7

# --- MID_LFT_#2 ---
$ws2 = $wb.Worksheets.Item(2)
$g2 = [double]"5.68432987514711e+23"
Add-DataRow $ws2 68 `
    "0x01,0x7c" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," `
    "0x01,0x4C" `
    "0x19" `
    380 `
    $g2 `
    332 `
    25

# --- MID_PLT_#1 ---
$ws3 = $wb.Worksheets.Item(3)
$g3 = [double]"5.68631262647113e+23"
Add-DataRow $ws3 68 `
    "0x00,0x6e" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," `
    "0x00,0x65" `
    "0x15" `
    110 `
    $g3 `
    101 `
    15

# --- MID_PLT_#2 ---
$ws4 = $wb.Worksheets.Item(4)
$g4 = [double]"5.68631262647113e+23"
Add-DataRow $ws4 68 `
    "0x00,0x82" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," `
    "0x00,0x7A" `
    "0x9" `
    130 `
    $g4 `
    122 `
    9
